$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 391
$ws.Range("I38").Value = 137
$ws.Range("J38").Value = 899
$ws.Range("K38").Value = 411
$ws.Range("L38").Value = 2697
$ws.Range("M38").Value = -39
$ws.Range("N38").Value = -3441
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9560
$ws.Range("H43").Value = 2119.75
$ws.Range("J43").Value = 2493
$ws.Range("L43").Value = 2493
$ws.Range("N43").Value = -2631
$ws.Range("H58").Value = 109.75
$ws.Range("J58").Value = 100
$ws.Range("L58").Value = 300
$ws.Range("N58").Value = -600
$ws.Range("H74").Value = 6079.6787
$ws.Range("J74").Value = 6549.6
$ws.Range("L74").Value = 6549.6
$ws.Range("N74").Value = -8421.6
$ws.Range("H76").Value = 111114744
$ws.Range("J76").Value = 4138.8
$ws.Range("L76").Value = 4138.8
$ws.Range("N76").Value = -4768.8
$ws.Range("H77").Value = 6079.6787
$ws.Range("J77").Value = 6549.6
$ws.Range("L77").Value = 32748
$ws.Range("N77").Value = -42108
$ws.Range("H79").Value = 111114744
$ws.Range("J79").Value = 4138.8
$ws.Range("L79").Value = 4138.8
$ws.Range("N79").Value = -6322.8
$ws.Range("H125").Value = 3875.5
$ws.Range("I125").Value = 3135
$ws.Range("J125").Value = 4122.3335
$ws.Range("K125").Value = 28215
$ws.Range("L125").Value = 37101.0015
$ws.Range("M125").Value = -25755
$ws.Range("N125").Value = -42021.0015
$ws.Range("H127").Value = 1202.3334
$ws.Range("I127").Value = 1053.5
$ws.Range("K127").Value = 3160.5
$ws.Range("M127").Value = 1799.5
$ws.Range("H129").Value = 1832.1666
$ws.Range("J129").Value = 3732.3333
$ws.Range("L129").Value = 11196.9999
$ws.Range("N129").Value = -21196.9999
$ws.Range("H137").Value = 771159.1
$ws.Range("I137").Value = 2001007
$ws.Range("K137").Value = 6003021
$ws.Range("M137").Value = -6000471
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 29000
$ws.Range("J24").Value = 29000
$ws.Range("L24").Value = 29000
$ws.Range("N24").Value = -29748
$ws.Range("H32").Value = 2657.866
$ws.Range("I32").Value = 1832.9615
$ws.Range("K32").Value = 1832.9615
$ws.Range("M32").Value = -1545.9615
$ws.Range("H45").Value = 4021.0715
$ws.Range("I45").Value = 3894.3333
$ws.Range("J45").Value = 4055.6365
$ws.Range("K45").Value = 3894.3333
$ws.Range("L45").Value = 4055.6365
$ws.Range("M45").Value = -3517.3333
$ws.Range("N45").Value = -4809.636500000001
$ws.Range("H95").Value = 45797.75
$ws.Range("J95").Value = 45797.75
$ws.Range("L95").Value = 45797.75
$ws.Range("N95").Value = -51289.75
$ws.Range("H97").Value = 1830.0646
$ws.Range("I97").Value = 652.73914
$ws.Range("J97").Value = 5214.875
$ws.Range("K97").Value = 652.73914
$ws.Range("L97").Value = 5214.875
$ws.Range("M97").Value = -156.73914
$ws.Range("N97").Value = -6206.875
$ws.Range("H100").Value = 29000
$ws.Range("J100").Value = 29000
$ws.Range("L100").Value = 29000
$ws.Range("N100").Value = -31164
$ws.Range("H102").Value = 4334.174
$ws.Range("I102").Value = 4711.2383
$ws.Range("K102").Value = 4711.2383
$ws.Range("M102").Value = -3089.2383
$ws.Range("H110").Value = 434.75
$ws.Range("I110").Value = 408
$ws.Range("K110").Value = 408
$ws.Range("M110").Value = 1637
$ws.Range("H122").Value = 2197.0625
$ws.Range("I122").Value = 1524.4445
$ws.Range("K122").Value = 4573.333500000001
$ws.Range("M122").Value = -2123.333500000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H86").Value = 6459.946
$ws.Range("I86").Value = 3600.12
$ws.Range("J86").Value = 12417.917
$ws.Range("K86").Value = 3600.12
$ws.Range("L86").Value = 12417.917
$ws.Range("M86").Value = -2477.12
$ws.Range("N86").Value = -14663.917
$ws.Range("H89").Value = 6459.946
$ws.Range("I89").Value = 3600.12
$ws.Range("J89").Value = 12417.917
$ws.Range("K89").Value = 18000.6
$ws.Range("L89").Value = 62089.585
$ws.Range("M89").Value = -12384.6
$ws.Range("N89").Value = -73321.58499999999
$ws.Range("H94").Value = 14671.786
$ws.Range("I94").Value = 450.41666
$ws.Range("J94").Value = 100000
$ws.Range("K94").Value = 450.41666
$ws.Range("L94").Value = 100000
$ws.Range("M94").Value = 0.5833400000000211
$ws.Range("N94").Value = -100902
$ws.Range("H97").Value = 33962.375
$ws.Range("I97").Value = 20204.666
$ws.Range("K97").Value = 20204.666
$ws.Range("M97").Value = -19213.666
$ws.Range("H99").Value = 2379.2
$ws.Range("I99").Value = 2379.2
$ws.Range("K99").Value = 2379.2
$ws.Range("M99").Value = -881.1999999999998
$ws.Range("H105").Value = 2639.1714
$ws.Range("I105").Value = 2049.476
$ws.Range("K105").Value = 2049.476
$ws.Range("M105").Value = -302.4760000000001
$ws.Range("H107").Value = 1682.3334
$ws.Range("I107").Value = 1438
$ws.Range("J107").Value = 2171
$ws.Range("K107").Value = 1438
$ws.Range("L107").Value = 2171
$ws.Range("M107").Value = 482
$ws.Range("N107").Value = -6011
$ws.Range("H141").Value = 61300
$ws.Range("J141").Value = 61950
$ws.Range("L141").Value = 61950
$ws.Range("N141").Value = -72310
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2164.516
$ws.Range("I31").Value = 2004.75
$ws.Range("J31").Value = 2265.4211
$ws.Range("K31").Value = 2004.75
$ws.Range("L31").Value = 2265.4211
$ws.Range("M31").Value = -1709.75
$ws.Range("N31").Value = -2855.4211
$ws.Range("H34").Value = 2164.516
$ws.Range("I34").Value = 2004.75
$ws.Range("J34").Value = 2265.4211
$ws.Range("K34").Value = 2004.75
$ws.Range("L34").Value = 2265.4211
$ws.Range("M34").Value = -1802.75
$ws.Range("N34").Value = -2669.4211
$ws.Range("H41").Value = 21424.334
$ws.Range("I41").Value = 7141.5
$ws.Range("K41").Value = 7141.5
$ws.Range("M41").Value = -6713.5
$ws.Range("H44").Value = 700
$ws.Range("J44").Value = 700
$ws.Range("L44").Value = 700
$ws.Range("N44").Value = -1584
$ws.Range("H55").Value = 2850
$ws.Range("J55").Value = 2850
$ws.Range("L55").Value = 2850
$ws.Range("N55").Value = -3480
$ws.Range("H58").Value = 2456.1875
$ws.Range("I58").Value = 1444.8
$ws.Range("K58").Value = 1444.8
$ws.Range("M58").Value = -1241.8
$ws.Range("H62").Value = 4021.3333
$ws.Range("I62").Value = 3082.25
$ws.Range("K62").Value = 3082.25
$ws.Range("M62").Value = -2458.25
$ws.Range("H65").Value = 4021.3333
$ws.Range("I65").Value = 3082.25
$ws.Range("K65").Value = 15411.25
$ws.Range("M65").Value = -12291.25
$ws.Range("H86").Value = 8281.714
$ws.Range("I86").Value = 7744.75
$ws.Range("J86").Value = 8997.666999999999
$ws.Range("K86").Value = 7744.75
$ws.Range("L86").Value = 8997.666999999999
$ws.Range("M86").Value = -6621.75
$ws.Range("N86").Value = -11243.667
$ws.Range("H89").Value = 8281.714
$ws.Range("I89").Value = 7744.75
$ws.Range("J89").Value = 8997.666999999999
$ws.Range("K89").Value = 38723.75
$ws.Range("L89").Value = 44988.335
$ws.Range("M89").Value = -33107.75
$ws.Range("N89").Value = -56220.335
$ws.Range("H134").Value = 3178.348
$ws.Range("I134").Value = 1297.9445
$ws.Range("K134").Value = 3893.8335
$ws.Range("M134").Value = -1358.8335
$ws.Range("H136").Value = 2456.1875
$ws.Range("I136").Value = 1444.8
$ws.Range("K136").Value = 4334.4
$ws.Range("M136").Value = -1784.4
$ws.Range("H140").Value = 105116.14
$ws.Range("J140").Value = 105116.14
$ws.Range("L140").Value = 105116.14
$ws.Range("N140").Value = -115476.14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 525
$ws.Range("I13").Value = 525
$ws.Range("K13").Value = 1575
$ws.Range("M13").Value = -1407
$ws.Range("H44").Value = 570
$ws.Range("I44").Value = 493.75
$ws.Range("J44").Value = 722.5
$ws.Range("K44").Value = 1481.25
$ws.Range("L44").Value = 2167.5
$ws.Range("M44").Value = -1083.25
$ws.Range("N44").Value = -2963.5
$ws.Range("H69").Value = 4473.3335
$ws.Range("I69").Value = 2400
$ws.Range("K69").Value = 7200
$ws.Range("M69").Value = -6389
$ws.Range("H72").Value = 4473.3335
$ws.Range("I72").Value = 2400
$ws.Range("K72").Value = 21600
$ws.Range("M72").Value = -17544
$ws.Range("H137").Value = 84284.86
$ws.Range("J137").Value = 11998.7
$ws.Range("L137").Value = 35996.10000000001
$ws.Range("N137").Value = -46196.10000000001
$ws.Range("H141").Value = 33242.133
$ws.Range("I141").Value = 9283
$ws.Range("J141").Value = 41954.547
$ws.Range("K141").Value = 27849
$ws.Range("L141").Value = 125863.641
$ws.Range("M141").Value = -22669
$ws.Range("N141").Value = -136223.641
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 11331.667
$ws.Range("I35").Value = 11000
$ws.Range("J35").Value = 12990
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 12990
$ws.Range("M35").Value = -10702
$ws.Range("N35").Value = -13586
$ws.Range("H38").Value = 29999
$ws.Range("J38").Value = 29999
$ws.Range("L38").Value = 29999
$ws.Range("N38").Value = -30925
$ws.Range("H40").Value = 31184.666
$ws.Range("J40").Value = 31184.666
$ws.Range("L40").Value = 31184.666
$ws.Range("N40").Value = -31486.666
$ws.Range("H80").Value = 3097.8333
$ws.Range("I80").Value = 2246.1667
$ws.Range("J80").Value = 3949.5
$ws.Range("K80").Value = 2246.1667
$ws.Range("L80").Value = 3949.5
$ws.Range("M80").Value = -1248.1667
$ws.Range("N80").Value = -5945.5
$ws.Range("H83").Value = 3097.8333
$ws.Range("I83").Value = 2246.1667
$ws.Range("J83").Value = 3949.5
$ws.Range("K83").Value = 11230.8335
$ws.Range("L83").Value = 19747.5
$ws.Range("M83").Value = -6238.833500000001
$ws.Range("N83").Value = -29731.5
$ws.Range("H97").Value = 905.2222
$ws.Range("I97").Value = 636.64
$ws.Range("K97").Value = 636.64
$ws.Range("M97").Value = -140.64
$ws.Range("H107").Value = 647.8570999999999
$ws.Range("I107").Value = 647.8570999999999
$ws.Range("K107").Value = 647.8570999999999
$ws.Range("M107").Value = 1272.1429
$ws.Range("H113").Value = 4119.75
$ws.Range("I113").Value = 1245
$ws.Range("J113").Value = 6994.5
$ws.Range("K113").Value = 1245
$ws.Range("L113").Value = 6994.5
$ws.Range("M113").Value = 925
$ws.Range("N113").Value = -11334.5
$ws.Range("H122").Value = 4763.64
$ws.Range("I122").Value = 3908.8235
$ws.Range("K122").Value = 11726.4705
$ws.Range("M122").Value = -9276.470499999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H16").Value = 417.5263
$ws.Range("I16").Value = 447.2143
$ws.Range("J16").Value = 334.4
$ws.Range("K16").Value = 447.2143
$ws.Range("L16").Value = 334.4
$ws.Range("M16").Value = -277.2143
$ws.Range("N16").Value = -674.4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H53").Value = 34999.5
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H82").Value = 1767.8572
$ws.Range("I82").Value = 887.5
$ws.Range("J82").Value = 2120
$ws.Range("K82").Value = 887.5
$ws.Range("L82").Value = 2120
$ws.Range("M82").Value = -526.5
$ws.Range("N82").Value = -2842
$ws.Range("H85").Value = 1767.8572
$ws.Range("I85").Value = 887.5
$ws.Range("J85").Value = 2120
$ws.Range("K85").Value = 887.5
$ws.Range("L85").Value = 2120
$ws.Range("M85").Value = 360.5
$ws.Range("N85").Value = -4616
$ws.Range("H88").Value = 99999
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 99999
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H93").Value = 1885.2354
$ws.Range("I93").Value = 1792.4546
$ws.Range("K93").Value = 1792.4546
$ws.Range("M93").Value = -544.4546
$ws.Range("H122").Value = 3509.6
$ws.Range("I122").Value = 3416
$ws.Range("J122").Value = 3650
$ws.Range("K122").Value = 10248
$ws.Range("L122").Value = 10950
$ws.Range("M122").Value = -7798
$ws.Range("N122").Value = -15850
$ws.Range("H132").Value = 2849.8225
$ws.Range("I132").Value = 2227.1943
$ws.Range("K132").Value = 6681.5829
$ws.Range("M132").Value = -4151.5829
$ws.Range("H136").Value = 2702.359
$ws.Range("J136").Value = 3056.5625
$ws.Range("L136").Value = 9169.6875
$ws.Range("N136").Value = -14269.6875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 33376.8
$ws.Range("I14").Value = 35498.5
$ws.Range("J14").Value = 24890
$ws.Range("K14").Value = 35498.5
$ws.Range("L14").Value = 24890
$ws.Range("M14").Value = -35330.5
$ws.Range("N14").Value = -25226
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H96").Value = 974.46155
$ws.Range("I96").Value = 1074.2
$ws.Range("J96").Value = 642
$ws.Range("K96").Value = 1074.2
$ws.Range("L96").Value = 642
$ws.Range("M96").Value = 298.8
$ws.Range("N96").Value = -3388
$ws.Range("H136").Value = 8890.286
$ws.Range("I136").Value = 7450.9585
$ws.Range("J136").Value = 17526.25
$ws.Range("K136").Value = 22352.8755
$ws.Range("L136").Value = 52578.75
$ws.Range("M136").Value = -19802.8755
$ws.Range("N136").Value = -57678.75
